$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(6, 4).Value = 13
$ws.Cells.Item(26, 4).Value = 12
$ws.Cells.Item(36, 4).Value = 10
$ws.Cells.Item(41, 4).Value = 27
$ws.Cells.Item(51, 4).Value = 21
$ws.Cells.Item(56, 4).Value = 7
$ws.Cells.Item(80, 4).Value = 6
$ws.Cells.Item(81, 4).Value = 26
$ws.Cells.Item(86, 4).Value = 6
$ws.Cells.Item(91, 4).Value = 14
$ws.Cells.Item(96, 4).Value = 10
$ws.Cells.Item(101, 4).Value = 14
$ws.Cells.Item(121, 4).Value = 7
$ws.Cells.Item(126, 4).Value = 16
$ws.Cells.Item(131, 4).Value = 3
$ws.Cells.Item(146, 4).Value = 2
$ws.Cells.Item(150, 4).Value = 25
$ws.Cells.Item(151, 4).Value = 20
$ws.Cells.Item(171, 4).Value = 4
$ws.Cells.Item(176, 4).Value = 14
$ws.Cells.Item(200, 4).Value = 27
$ws.Cells.Item(201, 4).Value = 21
$ws.Cells.Item(206, 4).Value = 6
$ws.Cells.Item(209, 4).Value = 6
$ws.Cells.Item(210, 4).Value = 5
$ws.Cells.Item(211, 4).Value = 10
$ws.Cells.Item(230, 4).Value = 7
$ws.Cells.Item(231, 4).Value = 17
$ws.Cells.Item(241, 4).Value = 17
$ws.Cells.Item(256, 4).Value = 25
$ws.Cells.Item(265, 4).Value = 16
$ws.Cells.Item(265, 5).Value = 4
$ws.Cells.Item(266, 4).Value = 21
$ws.Cells.Item(266, 5).Value = 2
$ws.Cells.Item(276, 4).Value = 14
$ws.Cells.Item(285, 4).Value = 9
$ws.Cells.Item(286, 4).Value = 6
$ws.Cells.Item(291, 4).Value = 3
$ws.Cells.Item(296, 4).Value = 7
$ws.Cells.Item(306, 4).Value = 2
$ws.Cells.Item(331, 4).Value = 12
$ws.Cells.Item(341, 4).Value = 4
$ws.Cells.Item(356, 4).Value = 14
$ws.Cells.Item(366, 4).Value = 6
$ws.Cells.Item(375, 4).Value = 10
$ws.Cells.Item(375, 5).Value = 2
$ws.Cells.Item(376, 4).Value = 11
$ws.Cells.Item(376, 5).Value = 4
$ws.Cells.Item(391, 4).Value = 3
$ws.Cells.Item(396, 4).Value = 25
$ws.Cells.Item(401, 4).Value = 22
$ws.Cells.Item(406, 4).Value = 25
$ws.Cells.Item(426, 4).Value = 10
$ws.Cells.Item(431, 4).Value = 15
$ws.Cells.Item(436, 4).Value = 17
$ws.Cells.Item(476, 4).Value = 21
$ws.Cells.Item(486, 4).Value = 2
$ws.Cells.Item(491, 4).Value = 26
$ws.Cells.Item(496, 4).Value = 7
$ws.Cells.Item(505, 4).Value = 6
$ws.Cells.Item(506, 4).Value = 25
$ws.Cells.Item(510, 4).Value = 22
$ws.Cells.Item(511, 4).Value = 11
$ws.Cells.Item(516, 4).Value = 16
$ws.Cells.Item(520, 4).Value = 22
$ws.Cells.Item(520, 5).Value = 4
$ws.Cells.Item(521, 4).Value = 1
$ws.Cells.Item(521, 5).Value = 3
$ws.Cells.Item(541, 4).Value = 17
$ws.Cells.Item(546, 4).Value = 2
$ws.Cells.Item(556, 4).Value = 13
$ws.Cells.Item(561, 4).Value = 16
$ws.Cells.Item(575, 4).Value = 22
$ws.Cells.Item(575, 5).Value = 4
$ws.Cells.Item(576, 4).Value = 3
$ws.Cells.Item(576, 5).Value = 5
$ws.Cells.Item(581, 4).Value = 17
$ws.Cells.Item(586, 4).Value = 6
$ws.Cells.Item(591, 4).Value = 26
$ws.Cells.Item(616, 4).Value = 6
$ws.Cells.Item(621, 4).Value = 13
$ws.Cells.Item(626, 4).Value = 9
$ws.Cells.Item(631, 4).Value = 7
$ws.Cells.Item(636, 4).Value = 17
$ws.Cells.Item(641, 4).Value = 25
$ws.Cells.Item(646, 4).Value = 12
$ws.Cells.Item(650, 4).Value = 17
$ws.Cells.Item(651, 4).Value = 6
$ws.Cells.Item(656, 4).Value = 22
$ws.Cells.Item(661, 4).Value = 1
$ws.Cells.Item(666, 4).Value = 8
$ws.Cells.Item(667, 4).Value = 25
$ws.Cells.Item(667, 5).Value = 2
$ws.Cells.Item(668, 4).Value = 8
$ws.Cells.Item(668, 5).Value = 1
$ws.Cells.Item(669, 4).Value = 19
$ws.Cells.Item(669, 5).Value = 3
$ws.Cells.Item(670, 4).Value = 12
$ws.Cells.Item(670, 5).Value = 5
$ws.Cells.Item(671, 4).Value = 11
$ws.Cells.Item(672, 4).Value = 16
$ws.Cells.Item(672, 5).Value = 5
$ws.Cells.Item(673, 4).Value = 5
$ws.Cells.Item(673, 5).Value = 1
$ws.Cells.Item(674, 4).Value = 15
$ws.Cells.Item(674, 5).Value = 3
$ws.Cells.Item(675, 4).Value = 26
$ws.Cells.Item(675, 5).Value = 2
$ws.Cells.Item(676, 4).Value = 1
$ws.Cells.Item(676, 5).Value = 4
